# svmvIncidentes.xlsx - "added new data for SNMV and OPM"
# Adds a 2024 data column (I), shifting the pre-existing "total" column to J
# with refreshed SUM formulas, and corrects a couple of 2023 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (Total de incidentes): 2023 figure correction ---
$ws.Range("H9").Value = 89

# --- New column I = year 2024 data (replaces the old "total" column) ---
$ws.Range("I1").Value = 2024
$ws.Range("I2").Value = 39
$ws.Range("I3").Value = 32
$ws.Range("I4").Value = 15
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 0
$ws.Range("I8").Value = 1

$ws.Range("I10").Value = 99

# --- Column J becomes the new "total" column, same look as old column I ---
# (copy formatting before I9 is switched to the F9:H9 number style below)
$ws.Range("I1:I10").Copy()
$ws.Range("J1:J10").PasteSpecial(-4122)
$ws.Range("J1").Value = "total"
$ws.Range("J2:J10").Formula = "=SUM(B2:I2)"

# I9 picks up the F9:H9 "total incidentes" number format (style index 2)
$ws.Range("H9").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Value = 95

# --- Cosmetic refresh matching the resaved workbook ---
$null = $ws.Range("G14").Select()
